$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: date, volumen, precio minimo, precio maximo, precio promedio ponderado, precio $/Kg
$ws.Range("D7").Value = 44355
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 6000
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 6000
$ws.Range("P7").Value = 375

# Row 8: date
$ws.Range("D8").Value = 44477
